$d = $word.ActiveDocument

# The <id>...</id> markers were previously split across three runs
# (<id>, the bare id text, </id>). Replace the whole tagged span in one
# shot so Word collapses it into a single run using the formatting of
# the first matched run (Courier New / color 7f6000 / sz 18), matching
# the newly downloaded tc/tcn/tl ids.
$d.Content.Find.Execute("<id>p048v_a1</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p048v_1</id>", 2)
$d.Content.Find.Execute("<id>p048v_a2</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p048v_2</id>", 2)
